# Autogenerated on Mon Feb 09 2015 03:30:35 GMT+0000 (Coordinated Universal Time)
# Insert a new "Number of employees" MSME size-classification table
# right after the existing MSME Participation table (row 15), pushing the
# Sector Distribution Details table (and everything below it) down by 6 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16-17 stay blank (same gap as before); insert 6 fresh rows at
# 21:26 (the first genuinely unused rows below the new table) so that the
# existing "Sector Distribution Details" block and all rows below shift
# down by 6 without leaving any stray placeholder rows behind.
$ws.Rows("21:26").Insert()

# Header row for the new table
$ws.Range("B18").Value = "Number of employees"
$ws.Range("C18").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D18").Value = "Turnover (local currency, unless noted otherwise)"

# Micro
$ws.Range("A19").Value = "Micro"
$ws.Range("B19").Value = "1-9"
$ws.Range("C19").Value = ""
$ws.Range("D19").Value = ""

# Small
$ws.Range("A20").Value = "Small"
$ws.Range("B20").Value = "10-49"
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = ""

# Medium
$ws.Range("A21").Value = "Medium"
$ws.Range("B21").Value = "50-300"
$ws.Range("C21").Value = ""
$ws.Range("D21").Value = ""

# Large
$ws.Range("A22").Value = "Large"
$ws.Range("B22").Value = ">300"
$ws.Range("C22").Value = ""
$ws.Range("D22").Value = ""
